# EnemyPartDrop.xlsx: re-order two pairs of data rows so the sheet matches
# huiji's expected JSON-array ordering.
#   - Row 8 ("Grineer Manic" ...) and Row 9 ("Ghoul Auger Alpha" ...) swap.
#   - Row 20 ("Zanuka Hunter" ...) and Row 21 ("Kuva Hyekka Master" ...) swap.
#
# We stage the existing row contents through scratch cells far below the
# used range (row 100+) using Copy/PasteSpecial (values), which preserves
# the shared-string cell typing and avoids introducing any new cell
# styles/number formats (unlike assigning `.Value` directly, which would
# coerce numeric-looking text like "0.5" into a real number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row 8 <-> row 9 -------------------------------------------------
$ws.Range("A8:J8").Copy()
$ws.Range("A100").PasteSpecial(-4163)

$ws.Range("A9:D9").Copy()
$ws.Range("A101").PasteSpecial(-4163)

$ws.Range("A8:J9").ClearContents()

$ws.Range("A101:D101").Copy()
$ws.Range("A8").PasteSpecial(-4163)

$ws.Range("A100:J100").Copy()
$ws.Range("A9").PasteSpecial(-4163)

$ws.Range("A100:J101").ClearContents()

# --- Swap row 20 <-> row 21 -----------------------------------------------
$ws.Range("A20:G20").Copy()
$ws.Range("A100").PasteSpecial(-4163)

$ws.Range("A21:D21").Copy()
$ws.Range("A101").PasteSpecial(-4163)

$ws.Range("A20:G21").ClearContents()

$ws.Range("A101:D101").Copy()
$ws.Range("A20").PasteSpecial(-4163)

$ws.Range("A100:G100").Copy()
$ws.Range("A21").PasteSpecial(-4163)

$ws.Range("A100:G101").ClearContents()

$excel.CutCopyMode = 0
